$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.867.56'
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').Value = '2.616.90'
$ws.Range('E3').Value = '  +2.87%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '517.31'
$ws.Range('E5').Value = '  +2.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.71'
$ws.Range('E6').Value = '  +1.28%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.565'
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = '2.639.34'
$ws.Range('E9').Value = '  +3.59%  '
$ws.Range('E10').Value = '  +3.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.104'
$ws.Range('E11').Value = '  +2.62%  '
$ws.Range('E12').Value = '  +1.40%  '
$ws.Range('D14').Value = '3.086.05'
$ws.Range('E14').Value = '  +3.31%  '
$ws.Range('D15').Value = '58.847.27'
$ws.Range('E15').Value = '  +0.78%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.94'
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('E17').Value = '  +1.23%  '
$ws.Range('D18').Value = '2.642.65'
$ws.Range('E18').Value = '  +3.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '349.51'
$ws.Range('E19').Value = '  +1.95%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.50'
$ws.Range('E20').Value = '  -0.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.29'
$ws.Range('E21').Value = '  +2.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.15'
$ws.Range('E22').Value = '  +3.53%  '
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.56'
$ws.Range('E24').Value = '  +1.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.418'
$ws.Range('E25').Value = '  +2.01%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.995'
$ws.Range('E26').Value = '  -0.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.162'
$ws.Range('E27').Value = '  +2.09%  '
$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').Value = '0.0₃0801'
$ws.Range('E28').Value = '  +2.07%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.07'
$ws.Range('E29').Value = '  +2.13%  '
$ws.Range('B30').Value = 'USDe'
$ws.Range('C30').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.998'
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.22'
$ws.Range('E31').Value = '  +7.16%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.94'
$ws.Range('E32').Value = '  +2.48%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.58'
$ws.Range('E33').Value = '  +2.96%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '149.86'
$ws.Range('E34').Value = '  +0.43%  '
$ws.Range('B35').Value = 'SuiNetwork'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.964'
$ws.Range('E35').Value = '  +6.85%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.99'
$ws.Range('E36').Value = '  +2.95%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.13'
$ws.Range('E37').Value = '  +2.46%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '36.67'
$ws.Range('E38').Value = '  +2.08%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.839'
$ws.Range('E39').Value = '  +2.75%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.68'
$ws.Range('E40').Value = '  +4.63%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.41'
$ws.Range('E41').Value = '  +1.61%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.995'
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '275.69'
$ws.Range('E43').Value = '  -1.56%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0983'
$ws.Range('E44').Value = '  -1.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.608'
$ws.Range('E45').Value = '  +1.21%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.55'
$ws.Range('E46').Value = '  +4.84%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0522'
$ws.Range('E47').Value = '  -1.68%  '
$ws.Range('B48').Value = 'WhiteBITCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.29'
$ws.Range('E48').Value = '  -0.09%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.71'
$ws.Range('E49').Value = '  +4.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0229'
$ws.Range('E50').Value = '  +1.03%  '
$ws.Range('D51').Value = '1.979.88'
$ws.Range('E51').Value = '  +4.66%  '
